# Refactor "Corrected" sheet (sheet2.xml) to drop the "adductName" column
# (column C), shifting the remaining columns (blank_1_404020,
# 072920_XXX2_1_TS1, 072920_XXX2_2_bra) left, and make the "Corrected"
# sheet the active/selected tab (it was "Original" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

# Delete column C ("adductName") -- shifts D:F left to C:E.
$ws.Range("C:C").Delete() | Out-Null

# Make the "Corrected" sheet the active sheet/tab.
$ws.Activate()
